$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3 ---
$ws.Range("D3").Value = 44482
$ws.Range("M3").Value = 240
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 11000
$ws.Range("P3").Value = 10500
$ws.Range("S3").Value = 5250

# --- Update row 4 ---
$ws.Range("D4").Value = 44454
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12500
$ws.Range("S4").Value = 6250

# --- Update row 5 ---
$ws.Range("D5").Value = 44475
$ws.Range("M5").Value = 240

# --- Add new row 6, reusing the same layout/style as row 5 ---
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44461
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11500
$ws.Range("Q6").Value = "$/bandeja 2 kilos"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 5750
$ws.Range("T6").Value = 2

# Match the date cell number format used by the other date cells (column D)
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
